$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.101.91'
$ws.Range('E2').Value = '  +0.90%  '

$ws.Range('D3').Value = '2.055.20'
$ws.Range('E3').Value = '  -3.17%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.654'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.07%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.85'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +15.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '61.90'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.71%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.377'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0785'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.08%  '

$ws.Range('E12').Value = '  +5.16%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.41%  '

$ws.Range('D14').Value = '2.355.37'
$ws.Range('E14').Value = '  -3.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.816'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.45%  '

$ws.Range('D17').Value = '2.056.06'
$ws.Range('E17').Value = '  -3.23%  '

$ws.Range('D18').Value = '37.049.96'
$ws.Range('E18').Value = '  +0.67%  '

$ws.Range('D19').Value = '0.0₃0909'
$ws.Range('E19').Value = '  +8.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.55%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.32%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.14%  '

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.66%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.51%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.97'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.37%  '

$ws.Range('E30').Value = '  -0.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.54'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.37%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0622'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.92%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.36%  '

$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0872'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.06%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.11%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.22%  '

$ws.Range('E39').Value = '  -0.07%  '

$ws.Range('E40').Value = '  +20.78%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.92%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0223'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.41%  '

$ws.Range('E43').Value = '  -4.85%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.76%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +46.55%  '

$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.77%  '

$ws.Range('B47').Value = 'Gas'
$ws.Range('C47').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -52.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.68%  '

$ws.Range('D49').Value = '1.295.38'
$ws.Range('E49').Value = '  -4.77%  '

$ws.Range('E50').Value = '  +2.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.08%  '
